# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit updates to the Leviathan_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 221.73914
$ws.Range("I9").Value = 256.85715
$ws.Range("J9").Value = 206.375
$ws.Range("K9").Value = 256.85715
$ws.Range("L9").Value = 206.375
$ws.Range("M9").Value = -87.85714999999999
$ws.Range("N9").Value = -544.375
$ws.Range("H12").Value = 462
$ws.Range("J12").Value = 698.5
$ws.Range("L12").Value = 698.5
$ws.Range("N12").Value = -1038.5
$ws.Range("H33").Value = 325
$ws.Range("I33").Value = 150
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 150
$ws.Range("L33").Value = 500
$ws.Range("M33").Value = 79
$ws.Range("N33").Value = -958
$ws.Range("H62").Value = 66960.75
$ws.Range("I62").Value = 75312.64
$ws.Range("J62").Value = 8497.5
$ws.Range("K62").Value = 75312.64
$ws.Range("L62").Value = 8497.5
$ws.Range("M62").Value = -74688.64
$ws.Range("N62").Value = -9745.5
$ws.Range("H64").Value = 3687.7576
$ws.Range("I64").Value = 3503.2222
$ws.Range("J64").Value = 4518.1665
$ws.Range("K64").Value = 3503.2222
$ws.Range("L64").Value = 4518.1665
$ws.Range("M64").Value = -3255.2222
$ws.Range("N64").Value = -5014.1665
$ws.Range("H65").Value = 66960.75
$ws.Range("I65").Value = 75312.64
$ws.Range("J65").Value = 8497.5
$ws.Range("K65").Value = 376563.2
$ws.Range("L65").Value = 42487.5
$ws.Range("M65").Value = -373443.2
$ws.Range("N65").Value = -48727.5
$ws.Range("H67").Value = 3687.7576
$ws.Range("I67").Value = 3503.2222
$ws.Range("J67").Value = 4518.1665
$ws.Range("K67").Value = 3503.2222
$ws.Range("L67").Value = 4518.1665
$ws.Range("M67").Value = -2645.2222
$ws.Range("N67").Value = -6234.1665
$ws.Range("H70").Value = 1062.8
$ws.Range("I70").Value = 882
$ws.Range("K70").Value = 2646
$ws.Range("M70").Value = -2376
$ws.Range("H73").Value = 1062.8
$ws.Range("I73").Value = 882
$ws.Range("K73").Value = 2646
$ws.Range("M73").Value = -1710
$ws.Range("H87").Value = 33354
$ws.Range("J87").Value = 33354
$ws.Range("L87").Value = 33354
$ws.Range("N87").Value = -35850
$ws.Range("H90").Value = 33354
$ws.Range("J90").Value = 33354
$ws.Range("L90").Value = 100062
$ws.Range("N90").Value = -112542
$ws.Range("H113").Value = 87240.75
$ws.Range("I113").Value = 500750
$ws.Range("J113").Value = 4538.9
$ws.Range("K113").Value = 500750
$ws.Range("L113").Value = 4538.9
$ws.Range("M113").Value = -497496
$ws.Range("N113").Value = -11046.9
$ws.Range("H116").Value = 5250
$ws.Range("I116").Value = 6000
$ws.Range("K116").Value = 6000
$ws.Range("M116").Value = -2558
$ws.Range("H132").Value = 3367.2307
$ws.Range("I132").Value = 1369.75
$ws.Range("J132").Value = 8451.727999999999
$ws.Range("K132").Value = 4109.25
$ws.Range("L132").Value = 25355.184
$ws.Range("M132").Value = -1579.25
$ws.Range("N132").Value = -30415.184
$ws.Range("H135").Value = 53930.79
$ws.Range("I135").Value = 1198.1666
$ws.Range("K135").Value = 10783.4994
$ws.Range("M135").Value = -8248.499400000001
$ws.Range("H137").Value = 3498.8635
$ws.Range("I137").Value = 1341.8438
$ws.Range("K137").Value = 4025.5314
$ws.Range("M137").Value = -1475.5314
$ws.Range("H138").Value = 2223.976
$ws.Range("I138").Value = 1247.0952
$ws.Range("J138").Value = 3200.8572
$ws.Range("K138").Value = 3741.2856
$ws.Range("L138").Value = 9602.571599999999
$ws.Range("M138").Value = 1398.7144
$ws.Range("N138").Value = -19882.5716
$ws.Range("H141").Value = 49626.527
$ws.Range("I141").Value = 49626.527
$ws.Range("K141").Value = 148879.581
$ws.Range("M141").Value = -143699.581

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2113.8096
$ws.Range("I74").Value = 968.62964
$ws.Range("K74").Value = 968.62964
$ws.Range("M74").Value = -94.62963999999999
$ws.Range("H77").Value = 2113.8096
$ws.Range("I77").Value = 968.62964
$ws.Range("K77").Value = 4843.1482
$ws.Range("M77").Value = -475.1481999999996
$ws.Range("H122").Value = 3990
$ws.Range("I122").Value = 2482
$ws.Range("J122").Value = 5875
$ws.Range("K122").Value = 7446
$ws.Range("L122").Value = 17625
$ws.Range("M122").Value = -4996
$ws.Range("N122").Value = -22525
$ws.Range("H132").Value = 14062.375
$ws.Range("I132").Value = 15507
$ws.Range("K132").Value = 46521
$ws.Range("M132").Value = -43991

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 2255000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 2255000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 2255000
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -2255278
$ws.Range("H31").Value = 2392.375
$ws.Range("I31").Value = 1822.8334
$ws.Range("J31").Value = 2961.9167
$ws.Range("K31").Value = 1822.8334
$ws.Range("L31").Value = 2961.9167
$ws.Range("M31").Value = -1527.8334
$ws.Range("N31").Value = -3551.9167
$ws.Range("H34").Value = 2392.375
$ws.Range("I34").Value = 1822.8334
$ws.Range("J34").Value = 2961.9167
$ws.Range("K34").Value = 1822.8334
$ws.Range("L34").Value = 2961.9167
$ws.Range("M34").Value = -1620.8334
$ws.Range("N34").Value = -3365.9167
$ws.Range("H86").Value = 7677.857
$ws.Range("I86").Value = 6936.5
$ws.Range("K86").Value = 6936.5
$ws.Range("M86").Value = -5813.5
$ws.Range("H89").Value = 7677.857
$ws.Range("I89").Value = 6936.5
$ws.Range("K89").Value = 34682.5
$ws.Range("M89").Value = -29066.5
$ws.Range("H107").Value = 1309.4
$ws.Range("I107").Value = 633.4
$ws.Range("K107").Value = 633.4
$ws.Range("M107").Value = 1286.6
$ws.Range("H132").Value = 4717.0557
$ws.Range("I132").Value = 5277.5713
$ws.Range("K132").Value = 15832.7139
$ws.Range("M132").Value = -13302.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2356.2856
$ws.Range("I39").Value = 996.4
$ws.Range("J39").Value = 2781.25
$ws.Range("K39").Value = 2989.2
$ws.Range("L39").Value = 8343.75
$ws.Range("M39").Value = -2695.2
$ws.Range("N39").Value = -8931.75
$ws.Range("H55").Value = 17859912
$ws.Range("I55").Value = 2951
$ws.Range("J55").Value = 41669196
$ws.Range("K55").Value = 8853
$ws.Range("L55").Value = 125007588
$ws.Range("M55").Value = -8676
$ws.Range("N55").Value = -125007942
$ws.Range("H115").Value = 2285.7144
$ws.Range("I115").Value = 1333.3334
$ws.Range("J115").Value = 3000
$ws.Range("K115").Value = 4000.0002
$ws.Range("L115").Value = 9000
$ws.Range("M115").Value = -2825.0002
$ws.Range("N115").Value = -11350
$ws.Range("H120").Value = 9382.25
$ws.Range("I120").Value = 9382.25
$ws.Range("K120").Value = 28146.75
$ws.Range("M120").Value = -23308.75
$ws.Range("H134").Value = 2416.3333
$ws.Range("I134").Value = 2416.3333
$ws.Range("K134").Value = 7248.999899999999
$ws.Range("M134").Value = -2178.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1868.3823
$ws.Range("I102").Value = 1644.4828
$ws.Range("K102").Value = 1644.4828
$ws.Range("M102").Value = -22.4828
$ws.Range("H122").Value = 716986.9
$ws.Range("I122").Value = 1252225
$ws.Range("K122").Value = 3756675
$ws.Range("M122").Value = -3754225
$ws.Range("H132").Value = 4562.5
$ws.Range("I132").Value = 4546.75
$ws.Range("J132").Value = 4578.25
$ws.Range("K132").Value = 13640.25
$ws.Range("L132").Value = 13734.75
$ws.Range("M132").Value = -11110.25
$ws.Range("N132").Value = -18794.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 26265.646
$ws.Range("I46").Value = 27819.75
$ws.Range("J46").Value = 1400
$ws.Range("K46").Value = 27819.75
$ws.Range("L46").Value = 1400
$ws.Range("M46").Value = -27631.75
$ws.Range("N46").Value = -1776
$ws.Range("H60").Value = 61000
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H122").Value = 3662.2307
$ws.Range("I122").Value = 3328.182
$ws.Range("K122").Value = 9984.545999999998
$ws.Range("M122").Value = -7534.545999999998
$ws.Range("H123").Value = 41000
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("H132").Value = 3079.56
$ws.Range("I132").Value = 2454.2727
$ws.Range("J132").Value = 7665
$ws.Range("K132").Value = 7362.8181
$ws.Range("L132").Value = 22995
$ws.Range("M132").Value = -4832.8181
$ws.Range("N132").Value = -28055

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 10322
$ws.Range("I4").Value = 40000
$ws.Range("K4").Value = 40000
$ws.Range("M4").Value = -39887

Write-Host "Applied 226 value updates and 3 clears."